$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# ------------------------------------------------------------------
# 1. Insert a brand-new blank row at row 78 -- this pushes the old
#    row 78 (and everything after it, through the last data row 142)
#    down by one, so the table needs one more row at the bottom.
# ------------------------------------------------------------------
$ws.Rows("78:78").Insert()

# Expand the table so it covers the newly added row at the bottom
# (A8:K142 -> A8:K143).
$tbl.Resize($ws.Range("A8:K143"))

# ------------------------------------------------------------------
# 2. The freshly inserted row 78 does not carry the table's normal
#    row formatting/formula, so rebuild it from row 79 (which now
#    holds what used to be row 78's content/format).
# ------------------------------------------------------------------
$ws.Range("A79:K79").Copy()
$ws.Range("A78:K78").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("G78").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Give A78 the "year separator" look (same format as A64 = "2023") and
# write the new "2024" year-label text into it.
$ws.Range("A64").Copy()
$ws.Range("A78").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A78").Value = "'2024"

# ------------------------------------------------------------------
# 3. Row 143 (old row 142, the specially-styled last table row) keeps
#    its own formatting after the shift, but its calculated-column
#    formula needs to be re-entered so it evaluates cleanly instead of
#    carrying a stale cached error from the resize.
# ------------------------------------------------------------------
$ws.Range("G143").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# ------------------------------------------------------------------
# 4. Leave-credit postings for the next few months (rows 74-76, which
#    sit above the inserted row so their row numbers don't move).
# ------------------------------------------------------------------
$ws.Range("C74").Value = 1.25

$ws.Range("B75").Value = "SL(1-0-0)"
$ws.Range("C75").Value = 1.25
$ws.Range("H75").Value = 1
$ws.Range("K74").Copy()
$ws.Range("K75").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K75").Value = 45202

$ws.Range("C76").Value = 1.25

# ------------------------------------------------------------------
# 5. Recalculate everything (BALANCE totals in E9/I9 pick up the new
#    EARNED postings automatically through the table formulas).
# ------------------------------------------------------------------
$excel.CalculateFull()

# Restore the cursor position recorded in the workbook.
$ws.Range("B68").Select()
